# Scheduled runner update: refresh market price data on the Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 292.57144
$ws.Range("I4").Value = 174.66667
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 174.66667
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -60.66667000000001
$ws.Range("N4").Value = -1228

$ws.Range("H19").Value = 11905396
$ws.Range("I19").Value = 41667010
$ws.Range("J19").Value = 751.6667
$ws.Range("K19").Value = 41667010
$ws.Range("L19").Value = 751.6667
$ws.Range("M19").Value = -41666835
$ws.Range("N19").Value = -1101.6667

$ws.Range("H43").Value = 1758.8572
$ws.Range("I43").Value = 750
$ws.Range("J43").Value = 1927
$ws.Range("K43").Value = 750
$ws.Range("L43").Value = 1927
$ws.Range("M43").Value = -681
$ws.Range("N43").Value = -2065

$ws.Range("H88").Value = 6048.5
$ws.Range("I88").Value = 6391.4287
$ws.Range("J88").Value = 5705.5713
$ws.Range("K88").Value = 6391.4287
$ws.Range("L88").Value = 5705.5713
$ws.Range("M88").Value = -5985.4287
$ws.Range("N88").Value = -6517.5713

$ws.Range("H91").Value = 6048.5
$ws.Range("I91").Value = 6391.4287
$ws.Range("J91").Value = 5705.5713
$ws.Range("K91").Value = 6391.4287
$ws.Range("L91").Value = 5705.5713
$ws.Range("M91").Value = -4987.4287
$ws.Range("N91").Value = -8513.5713

$ws.Range("H94").Value = 3085.2856
$ws.Range("J94").Value = 1798
$ws.Range("L94").Value = 1798
$ws.Range("N94").Value = -2700

$ws.Range("H116").Value = 2895.558
$ws.Range("I116").Value = 2344.4517
$ws.Range("J116").Value = 4319.25
$ws.Range("K116").Value = 2344.4517
$ws.Range("L116").Value = 4319.25
$ws.Range("M116").Value = 1097.5483
$ws.Range("N116").Value = -11203.25

$ws.Range("H132").Value = 7412988
$ws.Range("I132").Value = 10006160
$ws.Range("J132").Value = 3926.8572
$ws.Range("K132").Value = 30018480
$ws.Range("L132").Value = 11780.5716
$ws.Range("M132").Value = -30015950
$ws.Range("N132").Value = -16840.5716

$ws.Range("H135").Value = 546.70966
$ws.Range("I135").Value = 554.93335
$ws.Range("K135").Value = 4994.40015
$ws.Range("M135").Value = -2459.40015

$ws.Range("H138").Value = 3998.9697
$ws.Range("I138").Value = 2270.389
$ws.Range("J138").Value = 6073.2666
$ws.Range("K138").Value = 6811.167
$ws.Range("L138").Value = 18219.7998
$ws.Range("M138").Value = -1671.167
$ws.Range("N138").Value = -28499.7998

$ws.Range("H141").Value = 803424.75
$ws.Range("I141").Value = 1817.6818
$ws.Range("J141").Value = 2566960.5
$ws.Range("K141").Value = 5453.0454
$ws.Range("L141").Value = 7700881.5
$ws.Range("M141").Value = -273.0454
$ws.Range("N141").Value = -7711241.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 727.82355
$ws.Range("I74").Value = 533.7857
$ws.Range("K74").Value = 533.7857
$ws.Range("M74").Value = 340.2143

$ws.Range("H77").Value = 727.82355
$ws.Range("I77").Value = 533.7857
$ws.Range("K77").Value = 2668.9285
$ws.Range("M77").Value = 1699.0715

$ws.Range("H88").Value = 1500
$ws.Range("I88").Value = 1500
$ws.Range("K88").Value = 1500
$ws.Range("M88").Value = -1094

$ws.Range("H91").Value = 1500
$ws.Range("I91").Value = 1500
$ws.Range("K91").Value = 1500
$ws.Range("M91").Value = -96

$ws.Range("H102").Value = 2395.182
$ws.Range("I102").Value = 2429.7778
$ws.Range("J102").Value = 2239.5
$ws.Range("K102").Value = 2429.7778
$ws.Range("L102").Value = 2239.5
$ws.Range("M102").Value = -807.7777999999998
$ws.Range("N102").Value = -5483.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 17753.5
$ws.Range("I75").Value = 8107
$ws.Range("J75").Value = 27400
$ws.Range("K75").Value = 8107
$ws.Range("L75").Value = 27400
$ws.Range("M75").Value = -7171
$ws.Range("N75").Value = -29272

$ws.Range("H78").Value = 17753.5
$ws.Range("I78").Value = 8107
$ws.Range("J78").Value = 27400
$ws.Range("K78").Value = 24321
$ws.Range("L78").Value = 82200
$ws.Range("M78").Value = -19641
$ws.Range("N78").Value = -91560

$ws.Range("H86").Value = 1977.75
$ws.Range("I86").Value = 1286.3636
$ws.Range("J86").Value = 3498.8
$ws.Range("K86").Value = 1286.3636
$ws.Range("L86").Value = 3498.8
$ws.Range("M86").Value = -163.3635999999999
$ws.Range("N86").Value = -5744.8

$ws.Range("H89").Value = 1977.75
$ws.Range("I89").Value = 1286.3636
$ws.Range("J89").Value = 3498.8
$ws.Range("K89").Value = 6431.817999999999
$ws.Range("L89").Value = 17494
$ws.Range("M89").Value = -815.8179999999993
$ws.Range("N89").Value = -28726

$ws.Range("H99").Value = 1968.0588
$ws.Range("I99").Value = 1213.0834
$ws.Range("J99").Value = 3780
$ws.Range("K99").Value = 1213.0834
$ws.Range("L99").Value = 3780
$ws.Range("M99").Value = 284.9166
$ws.Range("N99").Value = -6776

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 48368.5
$ws.Range("I36").Value = 26684
$ws.Range("K36").Value = 26684
$ws.Range("M36").Value = -26296

$ws.Range("H40").Value = 48368.5
$ws.Range("I40").Value = 26684
$ws.Range("K40").Value = 26684
$ws.Range("M40").Value = -26524

$ws.Range("H62").Value = 3490.7222
$ws.Range("I62").Value = 2541.625
$ws.Range("J62").Value = 4250
$ws.Range("K62").Value = 2541.625
$ws.Range("L62").Value = 4250
$ws.Range("M62").Value = -1917.625
$ws.Range("N62").Value = -5498

$ws.Range("H65").Value = 3490.7222
$ws.Range("I65").Value = 2541.625
$ws.Range("J65").Value = 4250
$ws.Range("K65").Value = 12708.125
$ws.Range("L65").Value = 21250
$ws.Range("M65").Value = -9588.125
$ws.Range("N65").Value = -27490

$ws.Range("H99").Value = 2427.4546
$ws.Range("I99").Value = 1100
$ws.Range("J99").Value = 2722.4443
$ws.Range("K99").Value = 1100
$ws.Range("L99").Value = 2722.4443
$ws.Range("M99").Value = 398
$ws.Range("N99").Value = -5718.4443

$ws.Range("H126").Value = 2427.4546
$ws.Range("I126").Value = 1100
$ws.Range("J126").Value = 2722.4443
$ws.Range("K126").Value = 3300
$ws.Range("L126").Value = 8167.3329
$ws.Range("M126").Value = -830
$ws.Range("N126").Value = -13107.3329

$ws.Range("H132").Value = 4399.069
$ws.Range("I132").Value = 2833.25
$ws.Range("K132").Value = 8499.75
$ws.Range("M132").Value = -5969.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 7200
$ws.Range("I57").Value = 7200
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 7200
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -6380
$ws.Range("N57").ClearContents()

$ws.Range("H80").Value = 3761
$ws.Range("I80").Value = 3949.75
$ws.Range("K80").Value = 3949.75
$ws.Range("M80").Value = -2951.75

$ws.Range("H83").Value = 3761
$ws.Range("I83").Value = 3949.75
$ws.Range("K83").Value = 19748.75
$ws.Range("M83").Value = -14756.75

$ws.Range("H102").Value = 2675.8333
$ws.Range("I102").Value = 1901.375
$ws.Range("J102").Value = 4224.75
$ws.Range("K102").Value = 1901.375
$ws.Range("L102").Value = 4224.75
$ws.Range("M102").Value = -279.375
$ws.Range("N102").Value = -7468.75

$ws.Range("H126").Value = 2441.6775
$ws.Range("I126").Value = 1205.4117
$ws.Range("J126").Value = 3942.8572
$ws.Range("K126").Value = 3616.2351
$ws.Range("L126").Value = 11828.5716
$ws.Range("M126").Value = -1146.2351
$ws.Range("N126").Value = -16768.5716

$ws.Range("H139").Value = 29700
$ws.Range("J139").Value = 29700
$ws.Range("L139").Value = 29700
$ws.Range("N139").Value = -39980

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 590100.6
$ws.Range("I122").Value = 715665.9399999999
$ws.Range("J122").Value = 4129.3335
$ws.Range("K122").Value = 2146997.82
$ws.Range("L122").Value = 12388.0005
$ws.Range("M122").Value = -2144547.82
$ws.Range("N122").Value = -17288.0005
